$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 20 title text
$ws.Range("B20").Value = "Find the original typed string 1"

# Add new row 21
$ws.Range("A21").Value = 3333
$ws.Range("C21").Value = "Math/Dynamic Programming"
$ws.Range("B21").Value = "Find the original typed string 2"
$ws.Range("D21").Value = "Have an array of letter_groups count(same letter possiblility), compute totalCombination from each group, find number of invalid value for each group in an array, return the result - invalid value with MOD guard."

# Add new row 22
$ws.Range("A22").Value = 257
$ws.Range("B22").Value = "Binary Tree Paths"
$ws.Range("C22").Value = "Recursion"

# Update selection to D22 (matches the active cell / sqref left in the saved sheet view)
$ws.Range("D22").Select()
